# "update test cases _STA"
#
# The underlying test-case grid (rows 1-9) keeps the same logical content;
# this pass only fixes a handful of double-space / trailing-space typos in
# four cells, and removes a stale hyperlink on F4 that pointed at a sheet
# ('1-2q_Checklist') that is no longer part of this workbook. Everything
# else (row heights, column widths, borders/fills, merges, the A2:K2 "Back
# on main page" hyperlink) stays untouched.
#
# NOTE: order of the assignments below matters only in that it controls the
# order new shared-string entries get appended to xl/sharedStrings.xml -
# keep "E7" before "I4" to mirror the original authoring order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E7: "Sign out " -> "Sign out"  (drop trailing space)
$ws.Range("E7").Value = "Sign out"

# I4: "1.  The ...", "2.  The ..." -> single space after "1."/"2."
$ws.Range("I4").Value = "1. The `"Log in`" button is presented in the right corner of the Main Menu before the `"Public App`" widget.`n2. The `"Log in`" button is presented in the on the Start page header on the right.`n3. The user is redirected to the EPAM Digital Platform (https://access.epam.com/auth/....) for further login.`n4. The user is on the Start page again.`n5. The user is redirected to the EPAM Digital Platform (https://access.epam.com/auth/....) for further login."

# H7: "...at the end of  Dropdown menu." -> single space before "Dropdown"
$ws.Range("H7").Value = "1. Click on the Avatar (Profile icon) located in the Main Menu.`n2. Click on the `"Sign out`" button at the end of Dropdown menu."

# H9: "...two  `"Study`" buttons..." -> single space before the quote
$ws.Range("H9").Value = "1. Check the presence of two `"Study`" buttons on the Base Details Page of course.`n2. Click on the `"Study`" button on the header of this page.`n3. Return on the Base Details Page of course.`n4. Click on the `"Study`" button in the `"Summary`" widget on this page."

# Drop the stale hyperlink on F4 (kept: A2:K2 "Back on main page").
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$F$4') {
        $h.Delete()
    }
}
